# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" sheet, populated
#    with the per-fund holding breakdown for 2022-Q1.
# 2. Insert a new row at the top of the "总计" sheet's data with the 2022-Q1
#    aggregate totals (existing rows shift down by one).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: create the "2022-Q1" worksheet (placed immediately before "总计")
# ---------------------------------------------------------------------------

$template = $wb.Worksheets.Item("2021-Q1")
$totalSheetBefore = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheetBefore)
$newSheet.Name = "2022-Q1"
# Re-resolve "总计" by name: inserting a sheet shifts tab positions, and the
# handle obtained before the insert no longer points at the "总计" tab.
$totalSheet = $wb.Worksheets.Item("总计")

# Copy header formatting (bold + border) and the row-2 index-column style from
# an existing same-shaped sheet so the new sheet matches the established look.
$template.Range("B1:H1").Copy()
$newSheet.Range("B1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"

# Fund holding rows for 2022-Q1: code, name, scale, total position, position
# ratio, market value (亿元), position rank.
$fundRows = @(
    @("009909", "嘉实动力先锋混合A", "29.43", "90.99", "4.92", "1.4480", 5),
    @("005354", "富国沪港深行业精选灵活配置混合A", "37.00", "82.60", "2.37", "0.8769", 10),
    @("009869", "嘉实产业先锋混合A", "17.16", "89.91", "4.44", "0.7619", 6),
    @("009960", "银华多元机遇混合", "10.75", "87.53", "2.21", "0.2376", 9),
    @("009870", "嘉实产业先锋混合C", "3.71", "89.91", "4.44", "0.1647", 6),
    @("009910", "嘉实动力先锋混合C", "2.95", "90.99", "4.92", "0.1451", 5),
    @("040018", "华安香港精选股票(QDII)", "5.47", "88.46", "2.49", "0.1362", 7),
    @("011534", "万家民瑞祥明6个月持有期混合型证券投资基金A", "4.25", "23.07", "1.09", "0.0463", 2),
    @("011114", "富国沪港深行业精选灵活配置混合C", "1.68", "82.60", "2.37", "0.0398", 10),
    @("006786", "泰康中证港股通大消费主题指数A", "0.85", "80.77", "2.50", "0.0212", 9),
    @("006781", "汇丰晋信港股通精选股票", "0.67", "90.36", "3.15", "0.0211", 10),
    @("006787", "泰康中证港股通大消费主题指数C", "0.41", "80.77", "2.50", "0.0102", 9),
    @("011535", "万家民瑞祥明6个月持有期混合型证券投资基金C", "0.46", "23.07", "1.09", "0.0050", 2),
    @("009733", "创金合信港股通大消费精选股票A", "0.13", "82.28", "2.71", "0.0035", 10),
    @("009734", "创金合信港股通大消费精选股票C", "0.07", "82.28", "2.71", "0.0019", 10)
)

$r = 2
foreach ($row in $fundRows) {
    if ($r -gt 2) {
        $newSheet.Range("A2").Copy()
        $newSheet.Cells.Item($r, 1).PasteSpecial(-4122)
    }
    $newSheet.Cells.Item($r, 1).Value = $r - 2

    $dataRange = $newSheet.Range($newSheet.Cells.Item($r, 2), $newSheet.Cells.Item($r, 7))
    $dataRange.NumberFormat = "@"
    $newSheet.Cells.Item($r, 2).Value = $row[0]
    $newSheet.Cells.Item($r, 3).Value = $row[1]
    $newSheet.Cells.Item($r, 4).Value = $row[2]
    $newSheet.Cells.Item($r, 5).Value = $row[3]
    $newSheet.Cells.Item($r, 6).Value = $row[4]
    $newSheet.Cells.Item($r, 7).Value = $row[5]
    $dataRange.Style = "Normal"

    $newSheet.Cells.Item($r, 8).Value = $row[6]

    $r++
}

# ---------------------------------------------------------------------------
# Step 2: insert the 2022-Q1 summary row at the top of the "总计" sheet
# ---------------------------------------------------------------------------

$totalSheet.Rows.Item(2).Insert()

$totalSheet.Cells.Item(3, 1).Copy()
$totalSheet.Cells.Item(2, 1).PasteSpecial(-4122)
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 15
$totalSheet.Cells.Item(2, 4).Value = 3.92

# Column A is a simple 0-based row counter; re-stamp it for every data row
# now that a row was inserted at the top (Rows.Insert only shifts existing
# cell content down, it does not renumber it).
for ($i = 0; $i -lt 6; $i++) {
    $totalSheet.Cells.Item(2 + $i, 1).Value = $i
}
